$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), R (Origen), S (Precio $/Kg)
$rows = @(
    @{ Row = 2;  D = 44998; M = 20;  N = 2500; O = 2500; P = 2500; R = "Región de La Araucanía"; S = 2500 },
    @{ Row = 3;  D = 44616; M = 200; N = 3200; O = 3200; P = 3200; R = "Región de La Araucanía"; S = 3200 },
    @{ Row = 4;  D = 44567; M = 80;  N = 2400; O = 2400; P = 2400; R = "Región de La Araucanía"; S = 2400 },
    @{ Row = 5;  D = 44176; M = 20;  N = 3000; O = 3000; P = 3000; R = "Región de O'Higgins";    S = 3000 },
    @{ Row = 6;  D = 44592; M = 5;   N = 7500; O = 7500; P = 7500; R = "Región de La Araucanía"; S = 7500 },
    @{ Row = 7;  D = 44999; M = 25;  N = 2500; O = 2500; P = 2500; R = "Región de La Araucanía"; S = 2500 },
    @{ Row = 8;  D = 44574; M = 200; N = 3000; O = 3000; P = 3000; R = "Región de La Araucanía"; S = 3000 },
    @{ Row = 9;  D = 44175; M = 40;  N = 5000; O = 5000; P = 5000; R = "Provincia de Curicó";    S = 5000 },
    @{ Row = 10; D = 44214; M = 50;  N = 1800; O = 1800; P = 1800; R = "Región de La Araucanía"; S = 1800 },
    @{ Row = 11; D = 44551; M = 120; N = 4500; O = 4500; P = 4500; R = "Región de O'Higgins";    S = 4500 },
    @{ Row = 12; D = 44215; M = 65;  N = 2800; O = 2800; P = 2800; R = "Región de La Araucanía"; S = 2800 },
    @{ Row = 13; D = 44323; M = 20;  N = 3200; O = 3200; P = 3200; R = "Región de La Araucanía"; S = 3200 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Range("D$rowNum").Value = $r.D
    $ws.Range("M$rowNum").Value = $r.M
    $ws.Range("N$rowNum").Value = $r.N
    $ws.Range("O$rowNum").Value = $r.O
    $ws.Range("P$rowNum").Value = $r.P
    $ws.Range("R$rowNum").Value = $r.R
    $ws.Range("S$rowNum").Value = $r.S
}
